$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.229813
$ws.Range("H2").Value = 9.689439
$ws.Range("I2").Value = 0.6314649025465834
$ws.Range("J2").Value = 0.6663778861971119
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2654206666666667
$ws.Range("N2").Value = 0.796262
$ws.Range("O2").Value = 0.2863762486625487
$ws.Range("P2").Value = 0.2863762486625487
$ws.Range("Q2").Value = 0.8572591196686667
$ws.Range("R2").Value = 7.715332077018
$ws.Range("S2").Value = 0.1808365499533524
$ws.Range("T2").Value = 0.1908347992408077

$ws.Range("G3").Value = 3.229813
$ws.Range("H3").Value = 9.689439
$ws.Range("I3").Value = 0.6314649025465834
$ws.Range("J3").Value = 0.6663778861971119
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.6614043333333334
$ws.Range("N3").Value = 1.984213
$ws.Range("O3").Value = 0.7136237513374514
$ws.Range("P3").Value = 0.7136237513374514
$ws.Range("Q3").Value = 2.136212314056333
$ws.Range("R3").Value = 19.225910826507
$ws.Range("S3").Value = 0.450628352593231
$ws.Range("T3").Value = 0.4755430869563043

$ws.Range("G4").Value = 0.881166
$ws.Range("I4").Value = 0.1722779003977514
$ws.Range("J4").Value = 0.1818029515853594
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2654206666666667
$ws.Range("N4").Value = 0.796262
$ws.Range("O4").Value = 0.2863762486625487
$ws.Range("P4").Value = 0.2863762486625487
$ws.Range("Q4").Value = 0.233879667164
$ws.Range("R4").Value = 2.104917004476
$ws.Range("S4").Value = 0.04933629884336826
$ws.Range("T4").Value = 0.05206404727079417

$ws.Range("G5").Value = 0.881166
$ws.Range("I5").Value = 0.1722779003977514
$ws.Range("J5").Value = 0.1818029515853594
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.6614043333333334
$ws.Range("N5").Value = 1.984213
$ws.Range("O5").Value = 0.7136237513374514
$ws.Range("P5").Value = 0.7136237513374514
$ws.Range("Q5").Value = 0.5828070107860001
$ws.Range("R5").Value = 5.245263097074
$ws.Range("S5").Value = 0.1229416015543832
$ws.Range("T5").Value = 0.1297389043145652

$ws.Range("G6").Value = 0.1066033333333333
$ws.Range("H6").Value = 0.31981
$ws.Range("I6").Value = 0.02084215510138645
$ws.Range("J6").Value = 0.02199449439587765
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2654206666666667
$ws.Range("N6").Value = 0.796262
$ws.Range("O6").Value = 0.2863762486625487
$ws.Range("P6").Value = 0.2863762486625487
$ws.Range("Q6").Value = 0.02829472780222222
$ws.Range("R6").Value = 0.25465255022
$ws.Range("S6").Value = 0.005968698191978053
$ws.Range("T6").Value = 0.006298700796320891

$ws.Range("G7").Value = 0.1066033333333333
$ws.Range("H7").Value = 0.31981
$ws.Range("I7").Value = 0.02084215510138645
$ws.Range("J7").Value = 0.02199449439587765
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.6614043333333334
$ws.Range("N7").Value = 1.984213
$ws.Range("O7").Value = 0.7136237513374514
$ws.Range("P7").Value = 0.7136237513374514
$ws.Range("Q7").Value = 0.07050790661444445
$ws.Range("R7").Value = 0.63457115953
$ws.Range("S7").Value = 0.0148734569094084
$ws.Range("T7").Value = 0.01569579359955676

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8039255000000001
$ws.Range("H8").Value = 1.607851
$ws.Range("I8").Value = 0.1571765106872173
$ws.Range("J8").Value = 0.1105777486911175
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2654206666666667
$ws.Range("N8").Value = 0.796262
$ws.Range("O8").Value = 0.2863762486625487
$ws.Range("P8").Value = 0.2863762486625487
$ws.Range("Q8").Value = 0.2133784421603334
$ws.Range("R8").Value = 1.280270652962
$ws.Range("S8").Value = 0.04501161950847428
$ws.Range("T8").Value = 0.03166684085571227

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8039255000000001
$ws.Range("H9").Value = 1.607851
$ws.Range("I9").Value = 0.1571765106872173
$ws.Range("J9").Value = 0.1105777486911175
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.6614043333333334
$ws.Range("N9").Value = 1.984213
$ws.Range("O9").Value = 0.7136237513374514
$ws.Range("P9").Value = 0.7136237513374514
$ws.Range("Q9").Value = 0.5317198093771668
$ws.Range("R9").Value = 3.190318856263
$ws.Range("S9").Value = 0.112164891178743
$ws.Range("T9").Value = 0.0789109078354052

$ws.Range("G10").Value = 0.09328633333333332
$ws.Range("H10").Value = 0.279859
$ws.Range("I10").Value = 0.01823853126706141
$ws.Range("J10").Value = 0.01924691913053351
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2654206666666667
$ws.Range("N10").Value = 0.796262
$ws.Range("O10").Value = 0.2863762486625487
$ws.Range("P10").Value = 0.2863762486625487
$ws.Range("Q10").Value = 0.02476012078422222
$ws.Range("R10").Value = 0.222841087058
$ws.Range("S10").Value = 0.005223082165375648
$ws.Range("T10").Value = 0.00551186049891363

$ws.Range("G11").Value = 0.09328633333333332
$ws.Range("H11").Value = 0.279859
$ws.Range("I11").Value = 0.01823853126706141
$ws.Range("J11").Value = 0.01924691913053351
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.6614043333333334
$ws.Range("N11").Value = 1.984213
$ws.Range("O11").Value = 0.7136237513374514
$ws.Range("P11").Value = 0.7136237513374514
$ws.Range("Q11").Value = 0.06169998510744444
$ws.Range("R11").Value = 0.555299865967
$ws.Range("S11").Value = 0.01301544910168576
$ws.Range("T11").Value = 0.01373505863161988
